$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("WorkSheet")

$ws.Range("G7").Value = 0.3
$ws.Range("G10").Value = 76
$ws.Range("G11").Value = 88
$ws.Range("G12").Value = 91
$ws.Range("G17").Value = 1.7
$ws.Range("G18").Value = 4.6
$ws.Range("G21").Value = 0.098
$ws.Range("G26").Value = 1
$ws.Range("G28").Value = 2.4
$ws.Range("G29").Value = 0.975
